# Apply updated Team_PER_2013 data: re-order team codes in column B
# (fixing the PER-index/name mismatch bug) and replace column C
# values with the corrected per-minute PER figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "POR"
$ws.Range("C2").Value = 11.12142857142857
$ws.Range("B3").Value = "CLE"
$ws.Range("C3").Value = 12.32307692307692
$ws.Range("B4").Value = "DAL"
$ws.Range("C4").Value = 11.085
$ws.Range("B5").Value = "ATL"
$ws.Range("C5").Value = 14.17692307692308
$ws.Range("B6").Value = "OKC"
$ws.Range("C6").Value = 14.64615384615385
$ws.Range("B7").Value = "CHA"
$ws.Range("C7").Value = 11.95333333333333
$ws.Range("B8").Value = "WAS"
$ws.Range("C8").Value = 12.73846153846154
$ws.Range("B9").Value = "MIL"
$ws.Range("C9").Value = 13.41666666666667
$ws.Range("B10").Value = "LAC"
$ws.Range("C10").Value = 13.97857142857143
$ws.Range("B11").Value = "SAS"
$ws.Range("C11").Value = 14.72666666666666
$ws.Range("B12").Value = "DET"
$ws.Range("C12").Value = 13.48461538461538
$ws.Range("B13").Value = "ORL"
$ws.Range("C13").Value = 12.4
$ws.Range("B14").Value = "UTA"
$ws.Range("C14").Value = 13.59285714285714
$ws.Range("B15").Value = "MEM"
$ws.Range("C15").Value = 12.40909090909091
$ws.Range("B16").Value = "HOU"
$ws.Range("C16").Value = 12.91818181818182
$ws.Range("B17").Value = "DEN"
$ws.Range("C17").Value = 15.27333333333333
$ws.Range("B18").Value = "LAL"
$ws.Range("C18").Value = 10.6
$ws.Range("B19").Value = "GSW"
$ws.Range("C19").Value = 13
$ws.Range("B20").Value = "IND"
$ws.Range("C20").Value = 12.18666666666667
$ws.Range("B21").Value = "CHI"
$ws.Range("C21").Value = 12.68333333333334
$ws.Range("B22").Value = "PHI"
$ws.Range("C22").Value = 13.08461538461538
$ws.Range("B23").Value = "BOS"
$ws.Range("C23").Value = 11.05625
$ws.Range("B24").Value = "BRK"
$ws.Range("C24").Value = 12.225
$ws.Range("B25").Value = "TOR"
$ws.Range("C25").Value = 12.42307692307692
$ws.Range("B26").Value = "MIA"
$ws.Range("C26").Value = 14.06666666666667
$ws.Range("B27").Value = "SAC"
$ws.Range("C27").Value = 13.39090909090909
$ws.Range("B28").Value = "PHO"
$ws.Range("C28").Value = 12.46153846153846
$ws.Range("B29").Value = "NOH"
$ws.Range("C29").Value = 14.61428571428572
$ws.Range("B30").Value = "NYK"
$ws.Range("C30").Value = 13.02941176470588
$ws.Range("B31").Value = "MIN"
$ws.Range("C31").Value = 11.92777777777778
